$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "FFR" -> "C/A"
$ws.Range("B1").Value = "C/A"

# Row label: "FFR Lag" -> "C/A Lag"
$ws.Range("A2").Value = "C/A Lag"

# Coefficient values
$ws.Range("B2").Value = "-0.546**"
# "-0.014" looks numeric, so force it to be stored as text (leading
# apostrophe, same as typing it directly into Excel) rather than a number.
$ws.Range("C2").Value = "'-0.014"
$ws.Range("B3").Value = "1.055*"
$ws.Range("C3").Value = "0.922**"
